# Fixed :: Wrong value
# Reset the selector input cells (ASB level, AMP level, RB level) back to 0.
# These feed the INDEX/MATCH driven formulas (R7, S7, U7) which in turn
# cascade through the whole simulation table and the charts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 0
$ws.Range("R14").Value = 0

# Force a full recalculation so every dependent formula (K/L columns,
# the Z:AO helper table, and the chart caches) reflects the new inputs.
$excel.CalculateFullRebuild()

# Restore the cursor / scroll position that was captured in the saved file.
$ws.Activate()
$ws.Range("AS7").Select()
$excel.ActiveWindow.ScrollColumn = 16   # column P is the 16th column
$excel.ActiveWindow.ScrollRow = 2
